# CryCompanywiseStockReport_1.xlsx - stock quantity/value correction pass.
#
# The report is a flat, formula-free dump (every cell is a cached literal,
# not a live formula), so each "Qty" (col F) / "Value" (col G) pair and each
# "Sub Total:" (col B) row has to be poked with the already-recalculated
# numbers rather than relying on Excel to re-derive them. For a handful of
# duplicate-named items the two data rows simply swap their Item Code/MRP/
# Qty/Value (B/E/F/G) while the serial number (A), description (C) and Rate
# (D) stay put - so those are written as row-pairs below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 3M INDIA LTD ---------------------------------------------------------
# Row 9: 3M-Scotch Brite stainless steel Scrub (Combo Pack) - qty 146 -> 144
$ws.Range("F9").Value = 144
$ws.Range("G9").Value = 4258.08
# Row 10: Sub Total for the 3M INDIA LTD block
$ws.Range("B10").Value = 40780.8

# --- row 36: ZOFF Garlic powder 100 GMS - qty 8 -> 7 ---------------------
$ws.Range("F36").Value = 7
$ws.Range("G36").Value = 127.68
# Row 47: Sub Total
$ws.Range("B47").Value = 15143.33

# --- row 83: BHA-Vicks Vaporub 25 ml - qty 25 -> 22 -----------------------
$ws.Range("F83").Value = 22
$ws.Range("G83").Value = 1469.38

# --- row 93: GLT-7O CLOCK P II TWIN BLADE RAZOR - qty 72 -> 71 -----------
$ws.Range("F93").Value = 71
$ws.Range("G93").Value = 5607.58

# --- row 94: GLT-7O CLOCK PII TWIN BLADE 5S PACK - qty 149 -> 148 --------
$ws.Range("F94").Value = 148
$ws.Range("G94").Value = 20771.8
# Row 114: Sub Total
$ws.Range("B114").Value = 267020.12

# --- Brill International: several single-unit lines zeroed out -----------
# Row 119: BRILL-Imli (Tamarind) 200G - qty 1 -> 0
$ws.Range("F119").Value = 0
$ws.Range("G119").Value = 0
# Row 121: BRILL-Kalimirchi (Black Pepper) 100G - qty 2 -> 0
$ws.Range("F121").Value = 0
$ws.Range("G121").Value = 0
# Row 122: BRILL-Laung (Cloves) 100G - qty 1 -> 0
$ws.Range("F122").Value = 0
$ws.Range("G122").Value = 0
# Row 123: BRILL-Peanut Plain 200G - qty 4 -> 0
$ws.Range("F123").Value = 0
$ws.Range("G123").Value = 0
# Row 124: BRILL-Red Chilly (Lal Mirch) 200G - qty 1 -> 0
$ws.Range("F124").Value = 0
$ws.Range("G124").Value = 0
# Row 125: Sub Total
$ws.Range("B125").Value = 707.15

# --- rows 163/164: DAB-Real Activ Coconut Water Tetra 1000ml -------------
# the two lots swap Item Code / MRP / Qty / Value
$ws.Range("B163").Value = 57552
$ws.Range("E163").Value = 136.86
$ws.Range("F163").Value = -5
$ws.Range("G163").Value = -603.45
$ws.Range("B164").Value = 64329
$ws.Range("E164").Value = 128.32
$ws.Range("F164").Value = 3
$ws.Range("G164").Value = 362.07

# --- row 193: GHP-Glamic Cistem Toilet Cleaner Cube - qty 304 -> 303 -----
$ws.Range("F193").Value = 303
$ws.Range("G193").Value = 19634.4
# Row 200: Sub Total
$ws.Range("B200").Value = 48383.15

# --- row 229: HIM-ALMOND & ROSE SOAP 125G - qty 34 -> 31 -----------------
$ws.Range("F229").Value = 31
$ws.Range("G229").Value = 994.17

# --- rows 246/247: HIM-GENTLE BABY SOAP 75G swap --------------------------
$ws.Range("B246").Value = 64973
$ws.Range("E246").Value = 35.4
$ws.Range("F246").Value = 64
$ws.Range("G246").Value = 2131.2
$ws.Range("B247").Value = 48706
$ws.Range("E247").Value = 39.8
$ws.Range("F247").Value = -144
$ws.Range("G247").Value = -4795.2
# Row 274: Sub Total
$ws.Range("B274").Value = 89682.14

# --- row 284: HUL-Dove Dandruff Care Shampoo 180Ml - qty 30 -> 29 --------
$ws.Range("F284").Value = 29
$ws.Range("G284").Value = 4202.97

# --- row 285: HUL-Dove Dryness Care Shampoo 180Ml - qty 8 -> 7 -----------
$ws.Range("F285").Value = 7
$ws.Range("G285").Value = 922.8099999999999

# --- rows 292/293: HUL-Kissan nango jam 490g swap -------------------------
$ws.Range("B292").Value = 55373
$ws.Range("E292").Value = 163.62
$ws.Range("F292").Value = -94
$ws.Range("G292").Value = -13562.32
$ws.Range("B293").Value = 63520
$ws.Range("E293").Value = 153.4
$ws.Range("F293").Value = 73
$ws.Range("G293").Value = 10532.44

# --- rows 295/296: HUL-Kissan Pineapple Jam 500G swap ---------------------
$ws.Range("B295").Value = 63571
$ws.Range("E295").Value = 152.53
$ws.Range("F295").Value = 4
$ws.Range("G295").Value = 573.92
$ws.Range("B296").Value = 57802
$ws.Range("E296").Value = 162.71
$ws.Range("F296").Value = -79
$ws.Range("G296").Value = -11334.92

# --- rows 299/300: HUL-knorr schezwan 200g pch swap -----------------------
$ws.Range("B299").Value = 55356
$ws.Range("E299").Value = 54.04
$ws.Range("F299").Value = -158
$ws.Range("G299").Value = -7527.12
$ws.Range("B300").Value = 63510
$ws.Range("E300").Value = 50.66
$ws.Range("F300").Value = 145
$ws.Range("G300").Value = 6907.8

# --- row 321: HUL-Rin Bar 4X250G - qty 105 -> 104 ------------------------
$ws.Range("F321").Value = 104
$ws.Range("G321").Value = 6104.8

# --- row 328: HUL-Surf Excel Bar 250G - qty 824 -> 812 -------------------
$ws.Range("F328").Value = 812
$ws.Range("G328").Value = 17076.36

# --- row 333: HUL-Surf Exl Mtc Liq Fl 1 Ltr Cp - qty 719 -> 716 ----------
$ws.Range("F333").Value = 716
$ws.Range("G333").Value = 122672.28

# --- row 335: HUL-Surf Xl Mtk Pwd Poly Tl 1Kg - qty 9 -> 8 ---------------
$ws.Range("F335").Value = 8
$ws.Range("G335").Value = 1487.92
# Row 339: Sub Total
$ws.Range("B339").Value = 314211.87

# --- row 378: JYOTHY-MWP 1 KG - qty 8 -> 7 --------------------------------
$ws.Range("F378").Value = 7
$ws.Range("G378").Value = 519.12

# --- row 393: JYT - Ujala IDD Front load Detergent Liquid 2lt - qty 126 -> 125
$ws.Range("F393").Value = 125
$ws.Range("G393").Value = 26361.25
# Row 395: Sub Total
$ws.Range("B395").Value = 249551.95

# --- rows 420/421: KUS-Floor Wiper swap ------------------------------------
$ws.Range("B420").Value = 47097
$ws.Range("D420").Value = 112.28
$ws.Range("E420").Value = 134.16
$ws.Range("F420").Value = 15
$ws.Range("G420").Value = 1684.2
$ws.Range("B421").Value = 58047
$ws.Range("D421").Value = 105.54
$ws.Range("E421").Value = 126.1
$ws.Range("F421").Value = 42
$ws.Range("G421").Value = 4432.68

# --- rows 467/468: CRE-Butter cookies 100gm swap ---------------------------
$ws.Range("B467").Value = 53602
$ws.Range("E467").Value = 15.69
$ws.Range("F467").Value = -231
$ws.Range("G467").Value = -3037.65
$ws.Range("B468").Value = 65068
$ws.Range("E468").Value = 13.97
$ws.Range("F468").Value = 113
$ws.Range("G468").Value = 1485.95

# --- row 470: CRE-Cremfills 100gm - qty 139 -> 132 ------------------------
$ws.Range("F470").Value = 132
$ws.Range("G470").Value = 1690.92

# --- row 471: CRE-Cremica Chocochip Cookies (75+25Gm) - qty 71 -> 68 -----
$ws.Range("F471").Value = 68
$ws.Range("G471").Value = 1788.4

# --- rows 476/477: CRE-Cremica Golden Bytes Rich Butter 200Gm swap --------
$ws.Range("B476").Value = 64922
$ws.Range("E476").Value = 20.98
$ws.Range("F476").Value = 110
$ws.Range("G476").Value = 2170.3
$ws.Range("B477").Value = 45706
$ws.Range("E477").Value = 23.58
$ws.Range("F477").Value = -202
$ws.Range("G477").Value = -3985.46

# --- row 480: CRE-Cremica Honey Oatmeal Cookies 50+25 Gm - qty 198 -> 196
$ws.Range("F480").Value = 196
$ws.Range("G480").Value = 3179.12

# --- row 484: CRE-Cremica Melto 50Gm - qty 534 -> 528 ---------------------
$ws.Range("F484").Value = 528
$ws.Range("G484").Value = 3426.72

# --- row 485: CRE-Cremica Oatmeal Digestive 112.5 Gm - qty 203 -> 201 ----
$ws.Range("F485").Value = 201
$ws.Range("G485").Value = 2643.15

# --- row 488: CRE-Cremica Pista Almond Cookies (75+25Gm) - qty 107 -> 104
$ws.Range("F488").Value = 104
$ws.Range("G488").Value = 2735.2

# --- row 489: CRE-Cremica Sugar Crackers (100+20Gm) 120Gm - qty 69 -> 67 -
$ws.Range("F489").Value = 67
$ws.Range("G489").Value = 1100.81
# Row 492: Sub Total
$ws.Range("B492").Value = -2838.1

# --- row 559: UNB-McvitieS Cream Bourbon (100 Gms) - qty 106 -> 104 ------
$ws.Range("F559").Value = 104
$ws.Range("G559").Value = 2064.4
# Row 564: Sub Total
$ws.Range("B564").Value = 8121.19

# --- Pristine items: single-unit lines zeroed out -------------------------
# Row 570: PRI-P-01 Pristine Paper Nepkin - qty 2 -> 0
$ws.Range("F570").Value = 0
$ws.Range("G570").Value = 0
# Row 571: PRI-B 02 VIMAL Slim Gas Lighter (Stainless Steel) - qty 1 -> 0
$ws.Range("F571").Value = 0
$ws.Range("G571").Value = 0
# Row 573: PRI-B-31 Vimal Cloth clip Jumbo (stainless steel) Pack of 10 pcs - qty 3 -> 0
$ws.Range("F573").Value = 0
$ws.Range("G573").Value = 0
# Row 574: PRI-B-33 Vimal Fruit & Vegetable Juicer - qty 1 -> 0
$ws.Range("F574").Value = 0
$ws.Range("G574").Value = 0
# Row 579: Sub Total
$ws.Range("B579").Value = 9173.700000000001

# --- rows 603/604: Rasna Insta Orange 500g swap ---------------------------
$ws.Range("B603").Value = 64836
$ws.Range("E603").Value = 104.71
$ws.Range("F603").Value = 3
$ws.Range("G603").Value = 295.5
$ws.Range("B604").Value = 60031
$ws.Range("E604").Value = 111.69
$ws.Range("F604").Value = -5
$ws.Range("G604").Value = -492.5

# --- row 710: TCP-Chana Dal 1 kg - qty 47 -> 45 ---------------------------
$ws.Range("F710").Value = 45
$ws.Range("G710").Value = 3670.2

# --- row 712: TCP-Fine Besan 1 kg - qty 27 -> 26 --------------------------
$ws.Range("F712").Value = 26
$ws.Range("G712").Value = 2120.56

# --- row 713: TCP-kabooli chana 1 kg - qty 183 -> 181 ---------------------
$ws.Range("F713").Value = 181
$ws.Range("G713").Value = 25906.53

# --- row 714: TCP-Kala Chana 1 kg - qty 45 -> 43 --------------------------
$ws.Range("F714").Value = 43
$ws.Range("G714").Value = 3507.08

# --- rows 717/718: TCP-Rajma chitra 20x500g swap --------------------------
$ws.Range("B717").Value = 61428
$ws.Range("D717").Value = 69.16
$ws.Range("E717").Value = 73.52
$ws.Range("F717").Value = 1
$ws.Range("G717").Value = 69.16
$ws.Range("B718").Value = 63150
$ws.Range("D718").Value = 75.68000000000001
$ws.Range("E718").Value = 80.45
$ws.Range("F718").Value = 67
$ws.Range("G718").Value = 5070.56

# --- row 725: TCP-Toor Dal 1kg - qty 358 -> 357 ---------------------------
$ws.Range("F725").Value = 357
$ws.Range("G725").Value = 48198.57

# --- row 727: TCP-Urad Dal 1 kg - qty 280 -> 276 --------------------------
$ws.Range("F727").Value = 276
$ws.Range("G727").Value = 33315.96
# Row 728: Sub Total
$ws.Range("B728").Value = 146036.92

# --- row 732: Orgfeed Arhar Dal 1 Kg - qty 49 -> 48 -----------------------
$ws.Range("F732").Value = 48
$ws.Range("G732").Value = 7866.72

# --- row 736: Orgfeed Jeera 200 gm - qty 28 -> 27 -------------------------
$ws.Range("F736").Value = 27
$ws.Range("G736").Value = 4061.07

# --- rows 744/745: Shankys Tip Top Hing Jeera Peanut/Salted Peanut 200 Gm swap
$ws.Range("B744").Value = 65362
$ws.Range("F744").Value = 54
$ws.Range("G744").Value = 2206.98
$ws.Range("B745").Value = 65079
$ws.Range("F745").Value = 21
$ws.Range("G745").Value = 858.27

# --- row 752: Tip Top Maida 1 kg - qty 104 -> 103 -------------------------
$ws.Range("F752").Value = 103
$ws.Range("G752").Value = 5869.97

# --- row 753: Tip top MP Wheat Atta 5kg - qty 122 -> 120 ------------------
$ws.Range("F753").Value = 120
$ws.Range("G753").Value = 29016

# --- row 754: Tip Top Sooji 1 Kg - qty 116 -> 114 -------------------------
$ws.Range("F754").Value = 114
$ws.Range("G754").Value = 6496.86
# Row 755: Sub Total
$ws.Range("B755").Value = 83485.14999999999

# --- row 780: VVD Priyam Cold Pressed Groundnut Oil Pouch 1 Ltr - qty 3201 -> 3200
$ws.Range("F780").Value = 3200
$ws.Range("G780").Value = 521952

# --- row 781: VVD Priyam Groundnut Oil Bottle 1 Ltr - qty 63 -> 60 -------
$ws.Range("F781").Value = 60
$ws.Range("G781").Value = 10570.2

# --- row 783: VVD Pure Drop Cold Pressed Gingelly Oil Pouch 500Ml - qty 484 -> 483
$ws.Range("F783").Value = 483
$ws.Range("G783").Value = 69865.95

# --- row 786: VVD Veda Pancha Deepam Oil Pouch 900Ml - qty 181 -> 179 ----
$ws.Range("F786").Value = 179
$ws.Range("G786").Value = 23015.82
# Row 787: Sub Total
$ws.Range("B787").Value = 795194.9

# --- row 790: WIP-SOFTOUCH BLACK 800+200 ML - qty 109 -> 108 -------------
$ws.Range("F790").Value = 108
$ws.Range("G790").Value = 15769.08

# --- closing totals --------------------------------------------------------
# Row 804: Sub Total
$ws.Range("B804").Value = 79491.02
# Row 805: Sub Total (report-level)
$ws.Range("B805").Value = 3098417.73
# Row 806: Grand Total
$ws.Range("B806").Value = 3098417.73
